$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "C2" "16"
Set-TextValue "D2" "18"
Set-TextValue "E2" "1"
Set-TextValue "F2" "0"

Set-TextValue "C3" "20"
Set-TextValue "D3" "12"
Set-TextValue "E3" "2"
Set-TextValue "F3" "1"

Set-TextValue "C4" "1"
Set-TextValue "D4" "2"

Set-TextValue "C5" "0"

Set-TextValue "C7" "6"
Set-TextValue "D7" "4"

Set-TextValue "D8" "7"
